# Quarterly financials update for AVEO: insert a new reporting quarter
# (and the prior one) as columns D and E, shifting the existing quarterly
# history two columns to the right (old D:K -> F:M) on every financial
# statement block (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at D:E; this pushes the old D:K data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number formatting from the (now shifted) F:G columns into the
# freshly inserted D:E columns so the new cells pick up the same date /
# number styles as the rest of the table instead of the generic default.
# Each financial statement's data block is handled separately so we don't
# stamp formatting onto the blank header/spacer rows in between.
$dataBlocks = @(
    @{ First = 7;  Last = 35 },
    @{ First = 38; Last = 77 },
    @{ First = 80; Last = 102 }
)

foreach ($block in $dataBlocks) {
    $ws.Range("F$($block.First):G$($block.Last)").Copy()
    $ws.Range("D$($block.First):E$($block.Last)").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Populate the two new quarters of data.
$newData = @(
    @{ Row = 7; D = 43465; E = 43373 },
    @{ Row = 8; D = 1500; E = 2500 },
    @{ Row = 9; D = "NA"; E = "NA" },
    @{ Row = 10; D = "NA"; E = "NA" },
    @{ Row = 12; D = 5200; E = 5200 },
    @{ Row = 13; D = 0; E = 0 },
    @{ Row = 14; D = 0; E = 0 },
    @{ Row = 15; D = 0; E = 0 },
    @{ Row = 17; D = 7800; E = 7900 },
    @{ Row = 18; D = -6300; E = -5400 },
    @{ Row = 20; D = 28700; E = -16200 },
    @{ Row = 21; D = "NA"; E = "NA" },
    @{ Row = 22; D = 600; E = 600 },
    @{ Row = 23; D = 21800; E = -22200 },
    @{ Row = 24; D = 0; E = 0 },
    @{ Row = 25; D = 0; E = 0 },
    @{ Row = 26; D = 21800; E = -22200 },
    @{ Row = 27; D = 21800; E = -22200 },
    @{ Row = 28; D = 0; E = 0 },
    @{ Row = 29; D = 0; E = 0 },
    @{ Row = 30; D = 0; E = 0 },
    @{ Row = 31; D = 0; E = 0 },
    @{ Row = 32; D = -28700; E = 16200 },
    @{ Row = 33; D = 21800; E = -22200 },
    @{ Row = 34; D = 0; E = 0 },
    @{ Row = 35; D = 21800; E = -22200 },
    @{ Row = 38; D = 43465; E = 43373 },
    @{ Row = 41; D = 24400; E = 20400 },
    @{ Row = 42; D = 0; E = 0 },
    @{ Row = 43; D = 3000; E = 300 },
    @{ Row = 44; D = 0; E = 0 },
    @{ Row = 45; D = 500; E = 700 },
    @{ Row = 46; D = 27900; E = 21400 },
    @{ Row = 47; D = 0; E = 0 },
    @{ Row = 48; D = 0; E = 0 },
    @{ Row = 49; D = 0; E = 0 },
    @{ Row = 50; D = 0; E = 0 },
    @{ Row = 51; D = 0; E = 0 },
    @{ Row = 52; D = 0; E = 0 },
    @{ Row = 53; D = 0; E = 0 },
    @{ Row = 54; D = 27900; E = 21400 },
    @{ Row = 57; D = 3500; E = 3000 },
    @{ Row = 58; D = 3300; E = 4300 },
    @{ Row = 59; D = 11400; E = 12200 },
    @{ Row = 60; D = 18100; E = 19500 },
    @{ Row = 61; D = 15800; E = 14600 },
    @{ Row = 62; D = 21300; E = 47700 },
    @{ Row = 63; D = 0; E = 0 },
    @{ Row = 64; D = 0; E = 0 },
    @{ Row = 65; D = 0; E = 0 },
    @{ Row = 66; D = 55200; E = 81800 },
    @{ Row = 68; D = 0; E = 0 },
    @{ Row = 69; D = 0; E = 0 },
    @{ Row = 70; D = 0; E = 0 },
    @{ Row = 71; D = 0; E = 0 },
    @{ Row = 72; D = -595000; E = -616800 },
    @{ Row = 73; D = 0; E = 0 },
    @{ Row = 74; D = 0; E = 0 },
    @{ Row = 75; D = 0; E = 0 },
    @{ Row = 76; D = -27200; E = -60400 },
    @{ Row = 77; D = 0; E = 0 },
    @{ Row = 80; D = 43465; E = 43373 },
    @{ Row = 81; D = 21800; E = -22200 },
    @{ Row = 83; D = 0; E = 0 },
    @{ Row = 84; D = 0; E = 0 },
    @{ Row = 85; D = 0; E = 0 },
    @{ Row = 86; D = 0; E = 0 },
    @{ Row = 87; D = 0; E = 0 },
    @{ Row = 88; D = 0; E = 0 },
    @{ Row = 89; D = -6600; E = -2800 },
    @{ Row = 91; D = 0; E = 0 },
    @{ Row = 92; D = 0; E = 0 },
    @{ Row = 93; D = 0; E = 0 },
    @{ Row = 94; D = 0; E = 0 },
    @{ Row = 96; D = 0; E = 0 },
    @{ Row = 97; D = 0; E = 0 },
    @{ Row = 98; D = 0; E = 0 },
    @{ Row = 99; D = 0; E = 0 },
    @{ Row = 100; D = 10600; E = 5100 },
    @{ Row = 101; D = 0; E = 0 },
    @{ Row = 102; D = 4000; E = 2300 }
)

foreach ($entry in $newData) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D
    $ws.Cells.Item($entry.Row, 5).Value = $entry.E
}
